# Add a small caption textbox under the MTV-diagram picture on slide 10,
# crediting the image source (http://littlegreenriver.com/...).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# Position/size in points (exact EMU values from the target: off y=4561050 x=5979000,
# ext cy=331500 cx=2605499; 1 pt = 12700 EMU).
$left   = 5979000 / 12700
$top    = 4561050 / 12700
$width  = 2605499 / 12700
$height = 331500 / 12700

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "Shape 89"

$shp.Fill.Visible = $false
$shp.Line.Visible = $false

$shp.TextFrame.AutoSize = 0
$shp.TextFrame.MarginLeft = 91425 / 12700
$shp.TextFrame.MarginRight = 91425 / 12700
$shp.TextFrame.MarginTop = 91425 / 12700
$shp.TextFrame.MarginBottom = 91425 / 12700
$shp.TextFrame.VerticalAnchor = 1
$shp.TextFrame.HorizontalAnchor = 0

$shp.TextFrame.TextRange.Text = "(http://littlegreenriver.com/weblog/wp-content/uploads/mtv-diagram-730x1024.png)"
$shp.TextFrame.TextRange.ParagraphFormat.SpaceBefore = 0
$shp.TextFrame.TextRange.ParagraphFormat.Bullet.Visible = $false

$shp.TextFrame.TextRange.Font.Size = 8
$shp.TextFrame.TextRange.Font.Color.RGB = 13421772
